$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '25.326.32'
Set-TextValue 'E2' '  -0.36%  '
Set-TextValue 'D3' '1.659.74'
Set-TextValue 'E3' '  -0.65%  '
Set-TextValue 'D4' '0.9998'
Set-TextValue 'E4' '  -0.54%  '
Set-TextValue 'D5' '235.25'
Set-TextValue 'E5' '  -1.89%  '
Set-TextValue 'E6' '  -0.62%  '
Set-TextValue 'D7' '0.4777'
Set-TextValue 'E7' '  -1.10%  '
Set-TextValue 'D8' '0.2601'
Set-TextValue 'E8' '  -1.41%  '
Set-TextValue 'D9' '0.06142'
Set-TextValue 'E9' '  +2.04%  '
Set-TextValue 'D10' '0.07068'
Set-TextValue 'E10' '  -1.29%  '
Set-TextValue 'D11' '1.658.60'
Set-TextValue 'E11' '  -0.83%  '
Set-TextValue 'E12' '  +0.70%  '
Set-TextValue 'D13' '0.5894'
Set-TextValue 'E13' '  -6.12%  '
Set-TextValue 'D14' '4.379'
Set-TextValue 'E14' '  -5.93%  '
Set-TextValue 'D15' '74.27'
Set-TextValue 'E15' '  +0.82%  '
Set-TextValue 'E16' '  +0.05%  '
Set-TextValue 'D17' '1.001'
Set-TextValue 'E17' '  -0.49%  '
Set-TextValue 'D18' '25.312.73'
Set-TextValue 'E18' '  -0.43%  '
Set-TextValue 'D19' '0.000006738'
Set-TextValue 'E19' '  +1.34%  '
Set-TextValue 'D20' '11.40'
Set-TextValue 'E20' '  -1.58%  '
Set-TextValue 'D21' '1.872.20'
Set-TextValue 'D22' '4.432'
Set-TextValue 'E22' '  -0.86%  '
Set-TextValue 'D23' '8.626'
Set-TextValue 'E23' '  -0.10%  '
Set-TextValue 'D24' '5.319'
Set-TextValue 'E24' '  +0.14%  '
Set-TextValue 'D25' '133.51'
Set-TextValue 'E25' '  -0.45%  '
Set-TextValue 'E26' '  +1.24%  '
Set-TextValue 'D27' '1.401'
Set-TextValue 'E27' '  +1.47%  '
Set-TextValue 'D28' '104.02'
Set-TextValue 'E28' '  +1.17%  '
Set-TextValue 'D29' '1.686'
Set-TextValue 'E29' '  -2.28%  '
Set-TextValue 'D30' '3.966'
Set-TextValue 'E30' '  +2.09%  '
Set-TextValue 'D31' '3.618'
Set-TextValue 'E31' '  +1.59%  '
Set-TextValue 'D32' '0.07639'
Set-TextValue 'E32' '  -4.06%  '
Set-TextValue 'D33' '0.04358'
Set-TextValue 'E33' '  -6.02%  '
Set-TextValue 'D34' '1.000'
Set-TextValue 'E34' '  -0.59%  '
Set-TextValue 'D35' '2.603'
Set-TextValue 'E35' '  -1.67%  '
Set-TextValue 'D36' '0.6120'
Set-TextValue 'E36' '  +4.10%  '
Set-TextValue 'D37' '0.9435'
Set-TextValue 'E37' '  -1.16%  '
Set-TextValue 'D38' '2.607'
Set-TextValue 'E38' '  -1.51%  '
Set-TextValue 'D39' '0.8571'
Set-TextValue 'E39' '  +1.99%  '
Set-TextValue 'D40' '0.9998'
Set-TextValue 'E40' '  -0.69%  '
Set-TextValue 'E41' '  -3.89%  '
Set-TextValue 'D42' '1.826'
Set-TextValue 'E42' '  -2.79%  '
Set-TextValue 'D43' '97.75'
Set-TextValue 'E43' '  -1.46%  '
Set-TextValue 'D44' '0.3760'
Set-TextValue 'E44' '  -0.03%  '
Set-TextValue 'D45' '4.642'
Set-TextValue 'E45' '  -5.46%  '
Set-TextValue 'D46' '6.179'
Set-TextValue 'E46' '  +1.08%  '
Set-TextValue 'D47' '0.1107'
Set-TextValue 'E47' '  -3.73%  '
Set-TextValue 'D48' '0.05245'
Set-TextValue 'E48' '  +0.99%  '
Set-TextValue 'D49' '29.38'
Set-TextValue 'E49' '  -1.40%  '
Set-TextValue 'D50' '1.001'
Set-TextValue 'E50' '  -0.56%  '
Set-TextValue 'B51' 'EnergySwap'
Set-TextValue 'C51' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D51' '7.338'
Set-TextValue 'E51' '  -0.48%  '
